# Updates cryptos list data (coin names/links/prices/volume%) per the
# Wed Mar 13 16:42:49 UTC 2024 GitHub Actions refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that look like plain numbers (e.g. "1.00", "585.18") must be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# (e.g. "1.00" -> 1) and the formatted "D.DD"-style string is lost - exactly
# like typing such a value into a General-formatted cell in real Excel.
$textValueCells = [ordered]@{
    "D4" = "1.00"
    "D5" = "585.18"
    "D6" = "157.18"
    "D9" = "0.748"
    "D11" = "53.08"
    "D13" = "10.80"
    "D17" = "14.04"
    "D18" = "20.38"
    "D21" = "432.43"
    "D22" = "4.70"
    "D23" = "96.05"
    "D24" = "3.42"
    "D25" = "14.30"
    "D26" = "4.42"
    "D27" = "11.04"
    "D28" = "10.66"
    "D30" = "36.51"
    "D31" = "7.79"
    "D32" = "13.55"
    "D34" = "679.18"
    "D35" = "48.56"
    "D36" = "68.66"
    "D38" = "0.436"
    "D40" = "1.00"
    "D42" = "3.34"
    "D44" = "0.0485"
    "D45" = "10.71"
    "D47" = "2.65"
    "D50" = "3.42"
    "D51" = "2.15"
}
foreach ($ref in $textValueCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textValueCells[$ref]
}

# Remaining cells (coin names, links, multi-dot prices, and the
# "  +x.xx%  " volume strings) are not numeric-looking, so they round-trip
# as text without any extra formatting work.
$plainValueCells = [ordered]@{
    "D2" = "72.611.48"
    "E2" = "  +3.15%  "
    "D3" = "3.974.96"
    "E3" = "  +1.51%  "
    "E4" = "  -0.15%  "
    "E5" = "  +8.93%  "
    "E6" = "  +6.87%  "
    "E7" = "  -0.37%  "
    "E8" = "  -0.20%  "
    "E9" = "  +2.26%  "
    "E10" = "  +1.88%  "
    "E11" = "  +0.89%  "
    "E12" = "  +1.99%  "
    "E13" = "  +3.34%  "
    "D14" = "4.614.08"
    "E14" = "  +1.37%  "
    "D15" = "3.977.61"
    "E15" = "  +1.60%  "
    "E16" = "  +9.66%  "
    "E17" = "  +1.37%  "
    "E18" = "  +0.95%  "
    "E19" = "  -0.06%  "
    "D20" = "72.345.91"
    "E20" = "  +2.65%  "
    "E21" = "  +1.81%  "
    "E22" = "  +12.33%  "
    "E23" = "  +0.23%  "
    "E24" = "  -2.06%  "
    "E25" = "  +1.23%  "
    "E26" = "  +21.84%  "
    "E27" = "  -0.79%  "
    "E28" = "  +2.11%  "
    "E29" = "  +1.62%  "
    "E30" = "  +1.14%  "
    "E31" = "  +4.91%  "
    "E32" = "  +1.91%  "
    "E33" = "  +2.18%  "
    "B34" = "Bittensor"
    "C34" = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
    "E34" = "  +0.41%  "
    "B35" = "InjectiveProtocol"
    "C35" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "E35" = "  +1.77%  "
    "E36" = "  +5.23%  "
    "D37" = "0.0₃0883"
    "E37" = "  +8.59%  "
    "E38" = "  +2.21%  "
    "E39" = "  +1.41%  "
    "B40" = "Dai"
    "C40" = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
    "E40" = "  +0.00%  "
    "E41" = "  -1.62%  "
    "B42" = "ThetaToken"
    "C42" = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
    "E42" = "  -1.82%  "
    "E43" = "  +0.00%  "
    "E44" = "  +1.68%  "
    "E45" = "  +12.15%  "
    "E46" = "  +1.12%  "
    "E47" = "  -1.16%  "
    "E48" = "  +0.90%  "
    "E49" = "  +2.94%  "
    "E50" = "  +6.12%  "
    "E51" = "  +8.62%  "
}
foreach ($ref in $plainValueCells.Keys) {
    $ws.Range($ref).Value = $plainValueCells[$ref]
}
